$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row with the value "GRT-USD" in column A, row 41
$ws.Range("A41").Value = "GRT-USD"
